$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume updates (GitHub Actions data refresh)
# Cells whose new Price value is numeric-looking need the column
# kept as Text (matches the existing inlineStr cell type) so Excel
# does not silently coerce them into numbers.

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

$ws.Range("D2").Value = "60.983.05"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "3.382.10"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "570.88"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").Value = "141.33"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +2.19%  "
$ws.Range("D10").Value = "0.123"
$ws.Range("E10").Value = "  -0.91%  "
$ws.Range("E11").Value = "  -1.78%  "
$ws.Range("D12").Value = "3.960.22"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("D14").Value = "27.80"
$ws.Range("E14").Value = "  -1.31%  "
$ws.Range("D15").Value = "3.382.60"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "61.043.02"
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").Value = "6.11"
$ws.Range("E18").Value = "  -2.35%  "
$ws.Range("D19").Value = "13.60"
$ws.Range("E19").Value = "  -3.37%  "
$ws.Range("D20").Value = "8.93"
$ws.Range("E20").Value = "  -2.59%  "
$ws.Range("D21").Value = "383.43"
$ws.Range("E21").Value = "  -1.55%  "
$ws.Range("D22").Value = "75.52"
$ws.Range("E22").Value = "  +2.83%  "
$ws.Range("D23").Value = "0.551"
$ws.Range("E23").Value = "  -1.65%  "
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("E25").Value = "  -1.09%  "
$ws.Range("D26").Value = "3.523.37"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("E27").Value = "  +3.60%  "
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").Value = "7.21"
$ws.Range("E29").Value = "  -2.30%  "
$ws.Range("D30").Value = "7.96"
$ws.Range("E30").Value = "  -1.19%  "
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("E33").Value = "  -4.07%  "
$ws.Range("D34").Value = "23.21"
$ws.Range("E34").Value = "  -2.21%  "
$ws.Range("D35").Value = "6.94"
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("D36").Value = "165.88"
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("D37").Value = "3.415.62"
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("D38").Value = "4.97"
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("E39").Value = "  -2.44%  "
$ws.Range("D40").Value = "0.0765"
$ws.Range("E40").Value = "  -1.50%  "
$ws.Range("D41").Value = "26.71"
$ws.Range("E41").Value = "  -0.93%  "
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").Value = "0.778"
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("D44").Value = "4.36"
$ws.Range("E44").Value = "  -2.34%  "
$ws.Range("E45").Value = "  -2.52%  "
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("D47").Value = "2.450.27"
$ws.Range("E47").Value = "  -3.20%  "
$ws.Range("D48").Value = "22.90"
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("E49").Value = "  -2.73%  "
$ws.Range("E50").Value = "  +9.59%  "
